$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a numeric-looking string that must remain text
# (matches the original Price column formatting, e.g. "1.000", "0.04810").
# Setting NumberFormat to Text ("@") before assigning the value prevents Excel
# from auto-converting these into real numbers (which would lose formatting
# such as trailing zeros or dotted thousand separators).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D4").Value = '0.9999'
$ws.Range("D5").Value = '246.98'
$ws.Range("D7").Value = '0.4737'
$ws.Range("D8").Value = '0.2912'
$ws.Range("D9").Value = '0.06486'
$ws.Range("D10").Value = '22.17'
$ws.Range("D11").Value = '0.07715'
$ws.Range("D12").Value = '97.22'
$ws.Range("D13").Value = '0.7423'
$ws.Range("D15").Value = '5.153'
$ws.Range("D16").Value = '273.69'
$ws.Range("D18").Value = '13.39'
$ws.Range("D19").Value = '0.9994'
$ws.Range("D20").Value = '0.000007504'
$ws.Range("D22").Value = '1.0000'
$ws.Range("D23").Value = '5.255'
$ws.Range("D24").Value = '6.163'
$ws.Range("D25").Value = '9.282'
$ws.Range("D26").Value = '163.78'
$ws.Range("D27").Value = '18.74'
$ws.Range("D29").Value = '0.09996'
$ws.Range("D30").Value = '1.348'
$ws.Range("D31").Value = '1.508'
$ws.Range("D32").Value = '4.286'
$ws.Range("D33").Value = '4.099'
$ws.Range("D34").Value = '0.04806'
$ws.Range("D35").Value = '1.117'
$ws.Range("D36").Value = '0.6944'
$ws.Range("D38").Value = '0.01850'
$ws.Range("D39").Value = '2.743'
$ws.Range("D40").Value = '6.207'
$ws.Range("D41").Value = '73.19'
$ws.Range("D42").Value = '1.965'
$ws.Range("D43").Value = '0.4185'
$ws.Range("D44").Value = '0.9994'
$ws.Range("D45").Value = '0.8331'
$ws.Range("D46").Value = '102.47'
$ws.Range("D47").Value = '9.355'
$ws.Range("D48").Value = '35.37'
$ws.Range("D49").Value = '6.978'
$ws.Range("D50").Value = '923.63'
$ws.Range("D51").Value = '0.05640'

# Remaining cells (percentages and non-numeric-looking price strings) can be
# assigned directly since Excel will keep them as text already.
$ws.Range("D2").Value = '30.495.96'
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("D3").Value = '1.871.39'
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("E5").Value = '  +1.43%  '
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("E8").Value = '  +1.55%  '
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("E10").Value = '  +5.33%  '
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("E12").Value = '  +2.40%  '
$ws.Range("E13").Value = '  +5.48%  '
$ws.Range("D14").Value = '1.869.01'
$ws.Range("E14").Value = '  -0.01%  '
$ws.Range("E15").Value = '  +1.06%  '
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("D17").Value = '30.475.30'
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("E20").Value = '  -0.64%  '
$ws.Range("D21").Value = '2.114.93'
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("E23").Value = '  +0.94%  '
$ws.Range("E24").Value = '  +0.52%  '
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("E26").Value = '  -0.89%  '
$ws.Range("E27").Value = '  -0.91%  '
$ws.Range("E28").Value = '  +0.65%  '
$ws.Range("E29").Value = '  +1.53%  '
$ws.Range("E30").Value = '  -1.68%  '
$ws.Range("E32").Value = '  +1.05%  '
$ws.Range("E33").Value = '  +2.01%  '
$ws.Range("E34").Value = '  +1.63%  '
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("E36").Value = '  +0.46%  '
$ws.Range("E37").Value = '  +0.15%  '
$ws.Range("E38").Value = '  +0.53%  '
$ws.Range("E39").Value = '  +0.09%  '
$ws.Range("E40").Value = '  -1.94%  '
$ws.Range("E41").Value = '  +4.38%  '
$ws.Range("E42").Value = '  +3.81%  '
$ws.Range("E43").Value = '  +2.60%  '
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("E45").Value = '  -0.92%  '
$ws.Range("E46").Value = '  +0.60%  '
$ws.Range("E47").Value = '  +1.29%  '
$ws.Range("E49").Value = '  -1.15%  '
$ws.Range("E50").Value = '  +0.32%  '
$ws.Range("E51").Value = '  +1.34%  '
